$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.0345
$ws.Range("B9").Value = 8.415400000000005
$ws.Range("B18").Value = 4.635200000000005
$ws.Range("B20").Value = 5.807499999999999
